$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.472.26"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.569.89"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "'290.16"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "'0.3694"
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("D8").Value = "'49.98"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "'0.3382"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'1.146"
$ws.Range("E10").Value = "  +3.05%  "
$ws.Range("D11").Value = "'0.07534"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "'21.17"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").Value = "'6.013"
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("D15").Value = "'6.988"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").Value = "1.570.51"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "'0.00001121"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("D18").Value = "'90.26"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "'0.06768"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "'6.357"
$ws.Range("E21").Value = "  +4.00%  "
$ws.Range("D22").Value = "'16.34"
$ws.Range("D23").Value = "'12.15"
$ws.Range("E23").Value = "  +3.38%  "
$ws.Range("D24").Value = "22.467.34"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "'2.363"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "'2.655"
$ws.Range("E26").Value = "  +7.25%  "
$ws.Range("D27").Value = "'19.96"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "'149.59"
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").Value = "'5.037"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("D30").Value = "'124.56"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").Value = "1.746.20"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").Value = "'1.063"
$ws.Range("E32").Value = "  +9.73%  "
$ws.Range("D33").Value = "'6.219"
$ws.Range("E33").Value = "  +6.67%  "
$ws.Range("D34").Value = "'2.015"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("D35").Value = "'9.818"
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("D36").Value = "'0.08374"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "'0.02476"
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2294"
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.344"
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D40").Value = "'0.06509"
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("D41").Value = "'5.409"
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("D42").Value = "'11.25"
$ws.Range("E42").Value = "  +3.85%  "
$ws.Range("D43").Value = "'0.6222"
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("D44").Value = "'1.002"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "'14.02"
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("D46").Value = "'3.783"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'0.5860"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").Value = "'2.062"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("D49").Value = "'126.33"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").Value = "'1.234"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").Value = "'0.07294"
$ws.Range("E51").Value = "  +0.17%  "
